## Data_output. Corrigiendo funcion plot
## - Corrects montoOperado / volumenNominal / cantidadOperaciones for the
##   2022-05-16 and 2022-05-17 rows (92/93).
## - Appends the newly scraped daily quotes for 2024-05-14 .. 2024-05-28
##   (rows 578-588).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows 92 and 93: montoOperado (F), volumenNominal (G), cantidadOperaciones (H) ---
$ws.Cells.Item(92, 6).Value = 165022820.75
$ws.Cells.Item(92, 7).Value = 1000
$ws.Cells.Item(92, 8).Value = 944

$ws.Cells.Item(93, 6).Value = 257707074.1
$ws.Cells.Item(93, 7).Value = 2
$ws.Cells.Item(93, 8).Value = 1463

# --- Append new daily quote rows 578-588 ---
# Columns: fechaHora, ultimoPrecio, apertura, maximo, minimo, montoOperado, volumenNominal, cantidadOperaciones
$newRows = @(
    @("2024-05-14", 3816.25, 3615,    3825.9,   3615,    6438092387.9,      1713115, 5877),
    @("2024-05-15", 3989,    3873,    4000,     3850,    6971501378.9,      1752851, 7009),
    @("2024-05-16", 4003,    3990,    4022,     3900,    5073022575.6,      1278755, 6110),
    @("2024-05-17", 3889,    3995,    3998,     3880,    6641303031.25,     1685276, 6161),
    @("2024-05-20", 4088,    3880.1,  4100,     3830,    5813425285.95,     1446085, 6080),
    @("2024-05-21", 4200.9,  4110,    4205,     4016,    9302047470.450001, 2254635, 8724),
    @("2024-05-22", 4117.5,  4239.45, 4239.45,  4050,    9005214800.15,     2190266, 7406),
    @("2024-05-23", 3906,    4100,    4100,     3890,    9735366419.299999, 2447759, 8441),
    @("2024-05-24", 4004,    3905,    4055,     3890.8,  9337637896.85,     2350223, 6637),
    @("2024-05-27", 4077.95, 4099,    4099,     4040.15, 883192772.05,      240,     1929),
    @("2024-05-28", 4104,    4141,    4180,     4031,    6930788477.7,      1686361, 5395)
)

$startRow = 578
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Column A holds a plain "YYYY-MM-DD" text label (not a real Excel date),
    # matching the rest of the sheet -- force text so Excel doesn't
    # auto-convert it to a date serial, then drop back to the Normal style
    # so no stray number format lingers on the cell.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $data[0]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
